$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.901.63"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "3.157.58"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +32.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.370"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "3.156.72"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.95%  "
$ws.Range("E12").Value = "  +7.32%  "
$ws.Range("E13").Value = "  +7.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.01%  "
$ws.Range("D16").Value = "90.749.14"
$ws.Range("D17").Value = "3.740.76"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "3.129.30"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +15.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.11%  "
$ws.Range("D28").Value = "3.329.12"
$ws.Range("E28").Value = "  +3.84%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.73%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.209"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +51.21%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +20.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "520.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  +6.04%  "
$ws.Range("E37").Value = "  +6.99%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("E41").Value = "  +28.22%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.424"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.09%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +6.52%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.741"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +22.99%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.34%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "150.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.44%  "
